$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "Création de la map" line (row 44) ---------------------
# Mirror the formatting of an existing data row (A3:B3 -> "0%" red-fill
# data cell) onto the new row, then overwrite the values/text.
$ws.Range("A3:B3").Copy()
$ws.Range("A44:B44").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(44, 1).Value = "Création de la map"
$ws.Cells.Item(44, 2).Value = 0

# --- Update the active selection, as left by the author -----------------
[void]$ws.Range("C46").Select()

# --- Page setup: paper size / orientation were set on the sheet ---------
$ws.PageSetup.PaperSize = 9   # xlPaperA4
$ws.PageSetup.Orientation = 1 # xlPortrait
